$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# Remove the two "Old" helper rows (strComponentRootOld / strReportsDetailFilePathOld)
# which shifts every following row up and lets the engine garbage-collect the now
# orphaned shared strings (matches the target uniqueCount drop from 69 to 66).
$ws.Rows(8).Delete()
$ws.Rows(6).Delete()

# Rename the queue name value (drop the "_New" suffix)
$ws.Range("B3").Value = "LUX-01_LoadVendorOpenItems"

# Point the file-path settings at the new UNC share instead of the old local C: paths.
# Write B11 before B10 so the shared-string table append order matches the target file
# (LUX-01 Load Vendor Open Item Data path ends up before the Components path).
$ws.Range("B11").Value = "\\LRRBTUIPFSP100\Profiles\Uipath_26\Desktop\EssilorLuxottica Projects\LUX-01_Vendor Open Items Report\LUX - 01 - Load Vendor Open Item\Data"
$ws.Range("B10").Value = "\\LRRBTUIPFSP100\Profiles\Uipath_26\Desktop\EssilorLuxottica Projects\Components\"

# Move the hyperlink that used to live on B13 onto the new location of that same
# setting, B10, keeping the exact same target address.
$ws.Range("B13").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B10"), "file:///\\LRRBTUIPFSP100\Profiles\uipath_26\My%20Documents\UiPath\LUX%20-%2001%20-%20Load%20Vendor%20Open%20Item\Data")

# Restore the active selection recorded in the saved file.
$ws.Activate()
$ws.Range("B23").Select()
